# Updated cryptos list values (mirrors the scraped coinranking.com snapshot).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cell values, keyed by A1 reference.
$updates = [ordered]@{
    'D2' = '27.918.72'
    'E2' = '  +0.82%  '
    'D3' = '1.880.62'
    'E3' = '  +0.04%  '
    'D5' = '334.83'
    'E5' = '  +1.00%  '
    'D6' = '1.019'
    'E6' = '  +1.68%  '
    'D7' = '0.4648'
    'E7' = '  -1.62%  '
    'D9' = '46.86'
    'E9' = '  -2.13%  '
    'E10' = '  -1.80%  '
    'E11' = '  -2.08%  '
    'D12' = '21.50'
    'E12' = '  -1.60%  '
    'D13' = '1.895.20'
    'E13' = '  +0.85%  '
    'D14' = '5.905'
    'E14' = '  -1.05%  '
    'D15' = '7.056'
    'E15' = '  -1.71%  '
    'D16' = '1.022'
    'D17' = '0.06740'
    'E17' = '  +1.82%  '
    'D18' = '86.61'
    'E19' = '  -0.58%  '
    'E20' = '  -1.87%  '
    'E21' = '  +1.66%  '
    'D22' = '27.933.30'
    'E22' = '  +0.87%  '
    'D23' = '5.443'
    'E23' = '  -1.32%  '
    'E24' = '  -1.55%  '
    'E25' = '  +2.43%  '
    'D26' = '2.120.40'
    'E26' = '  +0.94%  '
    'D27' = '159.46'
    'E27' = '  +2.02%  '
    'E28' = '  -2.16%  '
    'D29' = '2.050'
    'E29' = '  -2.11%  '
    'D30' = '5.415'
    'E30' = '  -3.42%  '
    'D31' = '120.74'
    'E31' = '  -1.39%  '
    'D32' = '0.09451'
    'E32' = '  -1.30%  '
    'D33' = '0.9528'
    'E33' = '  -2.02%  '
    'D34' = '3.662'
    'E34' = '  +0.94%  '
    'D35' = '5.286'
    'E35' = '  -0.58%  '
    'D36' = '1.346'
    'E36' = '  -7.15%  '
    'D37' = '0.06077'
    'E37' = '  -0.72%  '
    'D38' = '0.02222'
    'E38' = '  -1.53%  '
    'E39' = '  -2.60%  '
    'D40' = '8.055'
    'E40' = '  -1.24%  '
    'D41' = '0.5854'
    'E41' = '  -2.75%  '
    'D42' = '0.1872'
    'E42' = '  -1.71%  '
    'E43' = '  -1.60%  '
    'D44' = '1.269'
    'E44' = '  +1.95%  '
    'D45' = '0.5601'
    'E45' = '  -1.82%  '
    'D46' = '12.02'
    'E46' = '  -1.68%  '
    'D47' = '3.386'
    'E47' = '  -0.67%  '
    'E48' = '  -1.73%  '
    'D49' = '0.06892'
    'E49' = '  +1.05%  '
    'D50' = '113.11'
    'E50' = '  +1.58%  '
    'B51' = 'BabyDogeCoin'
    'C51' = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
    'D51' = '0.00000000296'
    'E51' = '  -6.46%  '
}

# Some new values look like plain numbers (e.g. "334.83"); Excel would silently
# coerce those into numeric cells on assignment, so those specific cells are
# pre-formatted as Text to keep them as strings (matching the source data, which
# stores every "Price" entry as text).
$textForcedCells = @(
    'D5', 'D6', 'D7', 'D9', 'D12', 'D14', 'D15', 'D16', 'D17', 'D18', 'D23', 'D27', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D37', 'D38', 'D40', 'D41', 'D42', 'D44', 'D45', 'D46', 'D47', 'D49', 'D50', 'D51'
)
foreach ($cellRef in $textForcedCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

